$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 1.107673
$ws.Range("N2").Value = 3.323019
$ws.Range("O2").Value = 0.1862559687185926
$ws.Range("P2").Value = 0.1862559687185926
$ws.Range("Q2").Value = 223.1293640788147
$ws.Range("R2").Value = 2008.164276709332
$ws.Range("S2").Value = 0.09002682682150491
$ws.Range("T2").Value = 0.09002682682150491
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.2485388686743746
$ws.Range("P3").Value = 0.2485388686743746
$ws.Range("Q3").Value = 297.7425104693863
$ws.Range("R3").Value = 2679.682594224476
$ws.Range("S3").Value = 0.1201312679668618
$ws.Range("T3").Value = 0.1201312679668618
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 3.299541
$ws.Range("N4").Value = 9.898623000000001
$ws.Range("O4").Value = 0.55482006447906
$ws.Range("P4").Value = 0.55482006447906
$ws.Range("Q4").Value = 664.658689958116
$ws.Range("R4").Value = 5981.928209623045
$ws.Range("S4").Value = 0.2681722910980543
$ws.Range("T4").Value = 0.2681722910980544
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06176066666666667
$ws.Range("N5").Value = 0.185282
$ws.Range("O5").Value = 0.01038509812797287
$ws.Range("P5").Value = 0.01038509812797287
$ws.Range("Q5").Value = 12.44105280025511
$ws.Range("R5").Value = 111.969475202296
$ws.Range("S5").Value = 0.005019637422218192
$ws.Range("T5").Value = 0.005019637422218192
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 1.107673
$ws.Range("N6").Value = 3.323019
$ws.Range("O6").Value = 0.1862559687185926
$ws.Range("P6").Value = 0.1862559687185926
$ws.Range("Q6").Value = 72.46104857242067
$ws.Range("R6").Value = 652.1494371517859
$ws.Range("S6").Value = 0.02923612630756091
$ws.Range("T6").Value = 0.02923612630756092
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("O7").Value = 0.2485388686743746
$ws.Range("P7").Value = 0.2485388686743746
$ws.Range("S7").Value = 0.03901251491102935
$ws.Range("T7").Value = 0.03901251491102935
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("M8").Value = 3.299541
$ws.Range("N8").Value = 9.898623000000001
$ws.Range("O8").Value = 0.55482006447906
$ws.Range("P8").Value = 0.55482006447906
$ws.Range("Q8").Value = 215.847276829618
$ws.Range("R8").Value = 1942.625491466562
$ws.Range("S8").Value = 0.08708869624246131
$ws.Range("T8").Value = 0.08708869624246134
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06176066666666667
$ws.Range("N9").Value = 0.185282
$ws.Range("O9").Value = 0.01038509812797287
$ws.Range("P9").Value = 0.01038509812797287
$ws.Range("Q9").Value = 4.040220053389778
$ws.Range("R9").Value = 36.361980480508
$ws.Range("S9").Value = 0.001630122474327562
$ws.Range("T9").Value = 0.001630122474327563
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 1.107673
$ws.Range("N10").Value = 3.323019
$ws.Range("O10").Value = 0.1862559687185926
$ws.Range("P10").Value = 0.1862559687185926
$ws.Range("Q10").Value = 66.94204200421733
$ws.Range("R10").Value = 602.478378037956
$ws.Range("S10").Value = 0.02700935238834297
$ws.Range("T10").Value = 0.02700935238834297
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("O11").Value = 0.2485388686743746
$ws.Range("P11").Value = 0.2485388686743746
$ws.Range("Q11").Value = 89.3270669441898
$ws.Range("R11").Value = 803.9436024977081
$ws.Range("S11").Value = 0.03604112089620343
$ws.Range("T11").Value = 0.03604112089620343
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("M12").Value = 3.299541
$ws.Range("N12").Value = 9.898623000000001
$ws.Range("O12").Value = 0.55482006447906
$ws.Range("P12").Value = 0.55482006447906
$ws.Range("Q12").Value = 199.407236807828
$ws.Range("R12").Value = 1794.665131270452
$ws.Range("S12").Value = 0.08045557270853904
$ws.Range("T12").Value = 0.08045557270853905
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.06176066666666667
$ws.Range("N13").Value = 0.185282
$ws.Range("O13").Value = 0.01038509812797287
$ws.Range("P13").Value = 0.01038509812797287
$ws.Range("Q13").Value = 3.732496090640889
$ws.Range("R13").Value = 33.592464815768
$ws.Range("S13").Value = 0.001505963953024934
$ws.Range("T13").Value = 0.001505963953024934
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 1.107673
$ws.Range("N14").Value = 3.323019
$ws.Range("O14").Value = 0.1862559687185926
$ws.Range("P14").Value = 0.1862559687185926
$ws.Range("Q14").Value = 99.09856493454167
$ws.Range("R14").Value = 891.887084410875
$ws.Range("S14").Value = 0.03998366320118379
$ws.Range("T14").Value = 0.03998366320118379
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("O15").Value = 0.2485388686743746
$ws.Range("P15").Value = 0.2485388686743746
$ws.Range("Q15").Value = 132.2365419241806
$ws.Range("R15").Value = 1190.128877317625
$ws.Range("S15").Value = 0.05335396490028002
$ws.Range("T15").Value = 0.05335396490028002
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("M16").Value = 3.299541
$ws.Range("N16").Value = 9.898623000000001
$ws.Range("O16").Value = 0.55482006447906
$ws.Range("P16").Value = 0.55482006447906
$ws.Range("Q16").Value = 295.1952228163751
$ws.Range("R16").Value = 2656.757005347376
$ws.Range("S16").Value = 0.1191035044300052
$ws.Range("T16").Value = 0.1191035044300052
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.06176066666666667
$ws.Range("N17").Value = 0.185282
$ws.Range("O17").Value = 0.01038509812797287
$ws.Range("P17").Value = 0.01038509812797287
$ws.Range("Q17").Value = 5.525451497027778
$ws.Range("R17").Value = 49.72906347325001
$ws.Range("S17").Value = 0.00222937427840218
$ws.Range("T17").Value = 0.002229374278402181

Write-Output "Applied 182 cell updates"
